# Add a new "Greece" market tab, modeled on the existing "Croatia" tab
# (the last sheet in the workbook), with its own market name / item code.

$wb = $excel.ActiveWorkbook

# Leave the source sheet ("Croatia") with a "select-all" cursor state,
# mirroring how it was left after the new sheet was created and focus
# moved away from it.
$croatia = $wb.Worksheets.Item("Croatia")
$croatia.Activate()
$croatia.Cells.Select() | Out-Null

# Duplicate the Croatia sheet (same layout/styles/merged cells/page setup)
# and place the copy after it, as the new last tab.
$croatia.Copy($null, $wb.Worksheets.Item($wb.Worksheets.Count))
$greece = $wb.Worksheets.Item($wb.Worksheets.Count)
$greece.Name = "Greece"

# Fill in the market-specific values for Greece.
$greece.Range("B2").Value = "Greece Market"
$greece.Range("B4").Value = "NGC-4119/T3189"

# Make the new sheet the active tab/selection, matching how it was left
# after editing.
$greece.Activate()
$greece.Range("C19").Select() | Out-Null
